# Update comparison results for rows 4-8 (cornstover, sugarcane2g, oilcane2g, lactic, biorefinery)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comparison")

# Row 4 - cornstover
$ws.Range("C4").Value = 1.5335806045303599
$ws.Range("D4").Value = -12.7945415905519
$ws.Range("E4").Value = 1.5120574133501901
$ws.Range("F4").Value = -13.1766382676365

# Row 5 - sugarcane2g
$ws.Range("C5").Value = 2.1288498018910502
$ws.Range("E5").Value = 2.1292465237494902
$ws.Range("F5").Value = 0.72724579780821597

# Row 6 - oilcane2g
$ws.Range("C6").Value = 2.3535531988993199
$ws.Range("D6").Value = 1.76992480536927
$ws.Range("E6").Value = 2.3534974075271098
$ws.Range("F6").Value = 1.8590999877697201

# Row 7 - lactic
$ws.Range("C7").Value = 1.99972471216735
$ws.Range("D7").Value = 1.2669205590255099
$ws.Range("E7").Value = 1.9983229422348301
$ws.Range("F7").Value = 1.49180522952494

# Row 8 - biorefinery
$ws.Range("C8").Value = 1.4179202730054401
$ws.Range("D8").Value = 4.4994033158573004
$ws.Range("E8").Value = 1.3751680588605399
$ws.Range("F8").Value = 4.5054514989784797

# Update the active selection to F8 (matches diff for sheetView selection)
$ws.Range("F8").Select()

# Update workbook window position to match author's saved view
$excel.ActiveWindow.Left = 960
$excel.ActiveWindow.Top = 10180
